# "Treatment of literal queries, version 1"
# - Inserts a new "literal" column before the "stars" column (old O, new P).
# - Changes Level values for the two lemma rows (45, 48) to their more
#   specific nounlemma/verblemma labels.
# - Adds a new row (A051 / asta_lemma) describing the combined noun+verb
#   lemma filter.
# - Keeps the AutoFilter / _FilterDatabase range pinned at row 47 (as in the
#   source file) while covering the newly inserted column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at O - everything from O onward (stars, filter,
#    variants, unused1, unused2, Comments, ...) shifts one column right.
$ws.Columns("O:O").Insert()

# 2. Header + body for the new "literal" column.
$ws.Range("O1").Value = "literal"

# 3. The two "lemma" Level rows become more specific.
$ws.Range("D45").Value = "nounlemma"
$ws.Range("D48").Value = "verblemma"

# 4. New row 50: A051 / asta_lemma.
$ws.Range("A50").Value = "A051"
$ws.Range("D50").Value = "lemma"
$ws.Range("E50").Value = "lemma"
$ws.Range("E45").Copy()
$ws.Range("E50").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H50").Value = "yes"
$ws.Range("K50").Value = "asta_lemma"
$ws.Range("L50").Value = "yes"
$ws.Range("N50").Value = "CORE"
$ws.Range("O50").Value = "astalemmafunction"
$ws.Range("Q50").Value = "astalemmafilter"
$ws.Range("U50").Value = "noun en verb lemmas together"

# 5. Re-pin the AutoFilter / _FilterDatabase range to A1:W47 (matching the
#    original file, which deliberately stopped the filter range above the
#    last couple of data rows). Temporarily blank rows 48:50 so re-applying
#    the AutoFilter doesn't snap to the full contiguous data block, then
#    restore the saved values.
$lastRow = 50
$lastCol = 23
$saved = @{}
for ($r = 48; $r -le $lastRow; $r++) {
  for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $saved[[string]$r + "," + [string]$c] = $cell.Value2()
  }
}
$ws.Rows("48:" + $lastRow).ClearContents()
$ws.AutoFilterMode = $false
$ws.Range("A1:W47").AutoFilter()
foreach ($key in $saved.Keys) {
  $parts = $key.Split(",")
  $r = [int]$parts[0]
  $c = [int]$parts[1]
  $v = $saved[$key]
  if ($v -ne $null) {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

foreach ($n in $wb.Names) {
  if ($n.Name() -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$W`$47"
  }
}

# 6. Best-effort view state: selection lands on the new row/column, matching
#    where the author would have been after typing the new row.
$ws.Activate()
$ws.Range("O51").Select()
